# Auto-generated edit script applying targeted cell value updates
# per the commit diff for Garuda_Profits (market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2407.5293
$ws.Range("I40").Value = 4394.75
$ws.Range("J40").Value = 1796.0769
$ws.Range("K40").Value = 4394.75
$ws.Range("L40").Value = 1796.0769
$ws.Range("M40").Value = -4219.75
$ws.Range("N40").Value = -2146.0769

# Row 62
$ws.Range("H62").Value = 975
$ws.Range("I62").Value = 976.6
$ws.Range("J62").Value = 969.6667
$ws.Range("K62").Value = 976.6
$ws.Range("L62").Value = 969.6667
$ws.Range("M62").Value = -352.6
$ws.Range("N62").Value = -2217.6667

# Row 65
$ws.Range("H65").Value = 975
$ws.Range("I65").Value = 976.6
$ws.Range("J65").Value = 969.6667
$ws.Range("K65").Value = 4883
$ws.Range("L65").Value = 4848.3335
$ws.Range("M65").Value = -1763
$ws.Range("N65").Value = -11088.3335

# Row 86
$ws.Range("H86").Value = 52685332
$ws.Range("I86").Value = 101220.3
$ws.Range("J86").Value = 111112120
$ws.Range("K86").Value = 101220.3
$ws.Range("L86").Value = 111112120
$ws.Range("M86").Value = -100097.3
$ws.Range("N86").Value = -111114366

# Row 89
$ws.Range("H89").Value = 52685332
$ws.Range("I89").Value = 101220.3
$ws.Range("J89").Value = 111112120
$ws.Range("K89").Value = 506101.5
$ws.Range("L89").Value = 555560600
$ws.Range("M89").Value = -500485.5
$ws.Range("N89").Value = -555571832

# Row 135
$ws.Range("H135").Value = 3972.2104
$ws.Range("I135").Value = 689.931
$ws.Range("J135").Value = 14548.444
$ws.Range("K135").Value = 6209.379000000001
$ws.Range("L135").Value = 130935.996
$ws.Range("M135").Value = -3674.379000000001
$ws.Range("N135").Value = -136005.996

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 52632710
$ws.Range("I2").Value = 100000730
$ws.Range("J2").Value = 1580.3334
$ws.Range("K2").Value = 100000730
$ws.Range("L2").Value = 1580.3334
$ws.Range("M2").Value = -100000617
$ws.Range("N2").Value = -1806.3334

# Row 32
$ws.Range("H32").Value = 1204.89
$ws.Range("I32").Value = 1204.89
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1204.89
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -917.8900000000001
$ws.Range("N32").ClearContents()

# Row 80
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996

# Row 83
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984

# Row 88
$ws.Range("H88").Value = 1256964
$ws.Range("I88").Value = 1671133
$ws.Range("J88").Value = 14457
$ws.Range("K88").Value = 1671133
$ws.Range("L88").Value = 14457
$ws.Range("M88").Value = -1670727
$ws.Range("N88").Value = -15269

# Row 91
$ws.Range("H91").Value = 1256964
$ws.Range("I91").Value = 1671133
$ws.Range("J91").Value = 14457
$ws.Range("K91").Value = 1671133
$ws.Range("L91").Value = 14457
$ws.Range("M91").Value = -1669729
$ws.Range("N91").Value = -17265

# Row 116
$ws.Range("H116").Value = 52632710
$ws.Range("I116").Value = 100000730
$ws.Range("J116").Value = 1580.3334
$ws.Range("K116").Value = 100000730
$ws.Range("L116").Value = 1580.3334
$ws.Range("M116").Value = -99998436
$ws.Range("N116").Value = -6168.3334

# Row 140
$ws.Range("H140").Value = 48295.625
$ws.Range("J140").Value = 56611.816
$ws.Range("L140").Value = 56611.816
$ws.Range("N140").Value = -66971.81599999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 52632710
$ws.Range("I3").Value = 100000730
$ws.Range("J3").Value = 1580.3334
$ws.Range("K3").Value = 100000730
$ws.Range("L3").Value = 1580.3334
$ws.Range("M3").Value = -100000616
$ws.Range("N3").Value = -1808.3334

# Row 20
$ws.Range("H20").Value = 3877.8865
$ws.Range("I20").Value = 4362.724
$ws.Range("J20").Value = 2940.5334
$ws.Range("K20").Value = 4362.724
$ws.Range("L20").Value = 2940.5334
$ws.Range("M20").Value = -4115.724
$ws.Range("N20").Value = -3434.5334

# Row 35
$ws.Range("H35").Value = 29933.334
$ws.Range("J35").Value = 29933.334
$ws.Range("L35").Value = 29933.334
$ws.Range("N35").Value = -30553.334

# Row 86
$ws.Range("H86").Value = 2640.1155
$ws.Range("I86").Value = 2384.8823
$ws.Range("J86").Value = 3122.2222
$ws.Range("K86").Value = 2384.8823
$ws.Range("L86").Value = 3122.2222
$ws.Range("M86").Value = -1261.8823
$ws.Range("N86").Value = -5368.2222

# Row 89
$ws.Range("H89").Value = 2640.1155
$ws.Range("I89").Value = 2384.8823
$ws.Range("J89").Value = 3122.2222
$ws.Range("K89").Value = 11924.4115
$ws.Range("L89").Value = 15611.111
$ws.Range("M89").Value = -6308.411500000002
$ws.Range("N89").Value = -26843.111

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 334.4375
$ws.Range("I22").Value = 304.45456
$ws.Range("J22").Value = 400.4
$ws.Range("K22").Value = 304.45456
$ws.Range("L22").Value = 400.4
$ws.Range("M22").Value = 45.54543999999999
$ws.Range("N22").Value = -1100.4

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1818.14
$ws.Range("J131").Value = 847.2766
$ws.Range("L131").Value = 2541.8298
$ws.Range("N131").Value = -12621.8298

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 81.666664
$ws.Range("I2").Value = 51.833332
$ws.Range("J2").Value = 101.55556
$ws.Range("K2").Value = 51.833332
$ws.Range("L2").Value = 101.55556
$ws.Range("M2").Value = 61.166668
$ws.Range("N2").Value = -327.55556

# Row 52
$ws.Range("H52").Value = 11433.333
$ws.Range("J52").Value = 11433.333
$ws.Range("L52").Value = 11433.333
$ws.Range("N52").Value = -11951.333

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1258.3704
$ws.Range("I46").Value = 666
$ws.Range("J46").Value = 1332.4166
$ws.Range("K46").Value = 666
$ws.Range("L46").Value = 1332.4166
$ws.Range("M46").Value = -478
$ws.Range("N46").Value = -1708.4166

# Row 136
$ws.Range("H136").Value = 7247.6665
$ws.Range("I136").Value = 9328.643
$ws.Range("J136").Value = 3085.7144
$ws.Range("K136").Value = 27985.929
$ws.Range("L136").Value = 9257.143199999999
$ws.Range("M136").Value = -25435.929
$ws.Range("N136").Value = -14357.1432

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 933.2222
$ws.Range("I113").Value = 985.5714
$ws.Range("K113").Value = 2956.7142
$ws.Range("M113").Value = -786.7142000000003

# Row 136
$ws.Range("H136").Value = 1316.0968
$ws.Range("I136").Value = 1302.0975
$ws.Range("J136").Value = 1343.4286
$ws.Range("K136").Value = 3906.2925
$ws.Range("L136").Value = 4030.2858
$ws.Range("M136").Value = -1356.2925
$ws.Range("N136").Value = -9130.2858
